$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) values that look like plain decimal numbers need the cell
# pre-formatted as Text, otherwise Excel auto-converts the typed value
# into a floating point number instead of keeping the literal string.

$ws.Range("D2").Value = '57.871.43'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '3.123.84'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.49'
$ws.Range("E5").Value = '  +2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.03'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.493'
$ws.Range("E8").Value = '  +9.15%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  +3.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.139'
$ws.Range("E12").Value = '  +3.58%  '
$ws.Range("D13").Value = '3.663.43'
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000166'
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("D16").Value = '57.974.20'
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").Value = '3.121.09'
$ws.Range("E17").Value = '  +1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.12'
$ws.Range("E18").Value = '  +3.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.13'
$ws.Range("E21").Value = '  +7.60%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.40'
$ws.Range("E24").Value = '  +2.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.509'
$ws.Range("E25").Value = '  +1.95%  '
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  +3.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.17'
$ws.Range("E30").Value = '  +4.51%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.58'
$ws.Range("E32").Value = '  +3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.18'
$ws.Range("E33").Value = '  +4.42%  '
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '160.16'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E36").Value = '  +2.63%  '
$ws.Range("E37").Value = '  +4.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.67'
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("E39").Value = '  +5.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0671'
$ws.Range("D41").Value = '2.565.53'
$ws.Range("E41").Value = '  +7.46%  '
$ws.Range("E42").Value = '  +3.01%  '
$ws.Range("E43").Value = '  +4.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.697'
$ws.Range("E44").Value = '  +0.64%  '
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.982'
$ws.Range("E47").Value = '  +2.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.16'
$ws.Range("E48").Value = '  +3.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.89'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0937'
$ws.Range("E50").Value = '  +5.06%  '
$ws.Range("E51").Value = '  -1.30%  '
